# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.238.27"
$ws.Range("E2").Value = "  +5.11%  "
$ws.Range("D3").Value = "3.542.04"
$ws.Range("E3").Value = "  +5.73%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'189.60"
$ws.Range("E5").Value = "  +8.38%  "
$ws.Range("D6").Value = "'559.28"
$ws.Range("E6").Value = "  +4.90%  "
$ws.Range("D7").Value = "3.538.96"
$ws.Range("E7").Value = "  +5.88%  "
$ws.Range("E8").Value = "  +3.10%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("E11").Value = "  +13.10%  "
$ws.Range("D12").Value = "'54.89"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("E13").Value = "  +5.03%  "
$ws.Range("D14").Value = "'9.36"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "4.105.18"
$ws.Range("E15").Value = "  +5.86%  "
$ws.Range("D16").Value = "3.541.38"
$ws.Range("E16").Value = "  +5.86%  "
$ws.Range("D17").Value = "'0.122"
$ws.Range("E17").Value = "  +2.65%  "
$ws.Range("D18").Value = "67.219.60"
$ws.Range("E18").Value = "  +5.83%  "
$ws.Range("D19").Value = "'18.24"
$ws.Range("E19").Value = "  +3.61%  "
$ws.Range("D20").Value = "'12.06"
$ws.Range("E20").Value = "  +6.78%  "
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("D22").Value = "'429.89"
$ws.Range("E22").Value = "  +15.22%  "
$ws.Range("E23").Value = "  +9.00%  "
$ws.Range("D24").Value = "'85.07"
$ws.Range("E24").Value = "  +3.75%  "
$ws.Range("D25").Value = "'4.16"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "'11.04"
$ws.Range("E26").Value = "  -3.18%  "
$ws.Range("E27").Value = "  +6.75%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'6.15"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'12.33"
$ws.Range("E29").Value = "  +8.59%  "
$ws.Range("D30").Value = "'9.02"
$ws.Range("E30").Value = "  +8.16%  "
$ws.Range("D31").Value = "'30.47"
$ws.Range("E31").Value = "  +4.78%  "
$ws.Range("D32").Value = "'643.20"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").Value = "'6.61"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").Value = "'11.71"
$ws.Range("E34").Value = "  +3.60%  "
$ws.Range("E35").Value = "  +4.07%  "
$ws.Range("D36").Value = "'59.99"
$ws.Range("E36").Value = "  +2.38%  "
$ws.Range("D37").Value = "'38.33"
$ws.Range("E37").Value = "  +2.79%  "
$ws.Range("D38").Value = "0.0₃0811"
$ws.Range("E38").Value = "  +10.69%  "
$ws.Range("D39").Value = "'0.147"
$ws.Range("E39").Value = "  +16.66%  "
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("D41").Value = "'0.388"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").Value = "'3.42"
$ws.Range("E42").Value = "  +14.14%  "
$ws.Range("D43").Value = "3.119.40"
$ws.Range("E43").Value = "  +6.48%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("E45").Value = "  +2.39%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.36"
$ws.Range("E46").Value = "  +7.91%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.86"
$ws.Range("E47").Value = "  +8.60%  "
$ws.Range("D48").Value = "'0.0418"
$ws.Range("E48").Value = "  +4.25%  "
$ws.Range("E49").Value = "  +3.83%  "
$ws.Range("E50").Value = "  +4.78%  "
$ws.Range("D51").Value = "'143.30"
$ws.Range("E51").Value = "  +3.94%  "
